# Generate Report for Handback
#
# For both the "zh-cn" and "de-de" localization-status sheets:
#   - Row 2 (4345290d-...md) and Row 3 (b396911e-...md) are marked as
#     handed back (in sync with en-US) instead of "ready for handoff".
#   - The "Latest Target File" (E) / "Latest Handback File" (F) columns
#     are now populated with the source markdown file and the xlf file
#     that was handed back, each as a hyperlink (mirroring columns A/C).
#   - The "Latest Handback DateTime" (G) is stamped with the real handback
#     time (instead of the 0001-01-01 00:00:00 placeholder).
#   - The "Handoff Reason" (H) becomes "Include".

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

function Set-HandbackRow($ws, $row, $md, $mdUrl, $xlf, $xlfUrl, $handbackTime) {
    # Status column
    $ws.Range("B$row").Value = $statusText

    # Latest Target File (E) -- same file as "Source File Name" (A)
    $eCell = $ws.Range("E$row")
    $eCell.Value = $md
    $ws.Hyperlinks.Add($eCell, $mdUrl, "", "", $md)
    $eCell.Font.Name = "Calibri"
    $eCell.Font.Size = 11
    $eCell.Font.Underline = 2
    $eCell.Font.Color = 15570276

    # Latest Handback File (F) -- same file as "Latest Handoff File" (C)
    $fCell = $ws.Range("F$row")
    $fCell.Value = $xlf
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlf)
    $fCell.Font.Name = "Calibri"
    $fCell.Font.Size = 11
    $fCell.Font.Underline = 2
    $fCell.Font.Color = 15570276

    # Latest Handback DateTime (G)
    $ws.Range("G$row").Value = $handbackTime

    # Handoff Reason (H)
    $ws.Range("H$row").Value = "Include"
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow $wsZhCn 2 `
    "4345290d-7b98-49ac-89d3-937210843776.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/4345290d-7b98-49ac-89d3-937210843776.md" `
    "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd021bb156a4f150e364a328c690786dd9542736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf" `
    "2016-03-08 08:29:33"

Set-HandbackRow $wsZhCn 3 `
    "b396911e-8a1e-4350-bc5f-2848b741994d.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md" `
    "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd021bb156a4f150e364a328c690786dd9542736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf" `
    "2016-03-08 08:29:33"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow $wsDeDe 2 `
    "4345290d-7b98-49ac-89d3-937210843776.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/4345290d-7b98-49ac-89d3-937210843776.md" `
    "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5201e81b08a588b71f0b1ced1d2f3b1d44edf0e4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf" `
    "2016-03-08 08:29:40"

Set-HandbackRow $wsDeDe 3 `
    "b396911e-8a1e-4350-bc5f-2848b741994d.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md" `
    "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5201e81b08a588b71f0b1ced1d2f3b1d44edf0e4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf" `
    "2016-03-08 08:29:40"

Write-Output "Handback report generated."
